# Reposition/resize two texture-resource rectangles on the SLPK figure slide
# (slide 4), inside the "Group 8" group, to reflect the renamed texture
# resources / updated I3S 1.6 layout.
#
# NOTE: this runtime's GroupItems.Left/Top/Width/Height setters write the
# EMU value straight into the child shape's own <a:xfrm> (they do NOT apply
# the enclosing group's chOff/chExt rescale that the getters use), and they
# store the value as a single-precision (float32) point measurement before
# converting back to EMU (truncating, like real PowerPoint's COM layer).
# The literals below were chosen so that round(float32(pt) * 12700) lands
# exactly on the target EMU values from the target XML.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$g = $s.Shapes.Item("Group 8")

# "Rectangle 27" (Shared/ sharedResource.json.gz box):
#   off  x 6486853 -> 6358006   (y unchanged: 5034398)
#   ext cx 1675165 -> 1804012, cy 926182 -> 926181
$r27 = $g.GroupItems.Item("Rectangle 27")
$r27.Left   = 500.6304016113281
$r27.Width  = 142.04820251464844
$r27.Height = 72.92764282226562

# "Rectangle 19" (features/ 0.json.gz box):
#   off  x 4754192 -> 4754193   (y unchanged: 5041718)
#   ext cx 1654857 -> 1433403   (cy unchanged: 927017)
$r19 = $g.GroupItems.Item("Rectangle 19")
$r19.Left  = 374.3459167480469
$r19.Width = 112.86637878417969
